$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldUrl = "https://github.com/missysuperdoc/temp-doc-host/blob/main/Blank.docx"
$newUrl = "https://github.com/missysuperdoc/temp-doc-host/blob/main/Starting%20files/Blank.docx"

$cells = "F2", "F3", "F4", "F5", "F6"

foreach ($cellRef in $cells) {
    $range = $ws.Range($cellRef)
    $current = [string]$range.Value2
    $range.Value2 = $current.Replace($oldUrl, $newUrl)
}
